$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("enc_route")

$ws.Range("A2").Value2 = "D"
$ws.Range("B2").Value2 = 7
$ws.Range("A3").Value2 = "V"
$ws.Range("B3").Value2 = 3
$ws.Range("A4").Value2 = "o"
$ws.Range("B4").Value2 = 68
$ws.Range("A5").Value2 = "."
$ws.Range("B5").Value2 = 35
$ws.Range("A6").Value2 = "c"
$ws.Range("B6").Value2 = 66
$ws.Range("A7").Value2 = "v"
$ws.Range("B7").Value2 = 22
$ws.Range("A8").Value2 = "x"
$ws.Range("B8").Value2 = 3
$ws.Range("A9").Value2 = "`n"
$ws.Range("B9").Value2 = 2
$ws.Range("A10").Value2 = "l"
$ws.Range("B10").Value2 = 81
$ws.Range("A11").Value2 = ","
$ws.Range("B11").Value2 = 19
$ws.Range("A12").Value2 = "g"
$ws.Range("B12").Value2 = 18
$ws.Range("A13").Value2 = "h"
$ws.Range("B13").Value2 = 15
$ws.Range("A14").Value2 = "I"
$ws.Range("B14").Value2 = 5
$ws.Range("A15").Value2 = "i"
$ws.Range("B15").Value2 = 130
$ws.Range("A16").Value2 = "n"
$ws.Range("B16").Value2 = 83
$ws.Range("A17").Value2 = "P"
$ws.Range("B17").Value2 = 5
$ws.Range("A18").Value2 = "j"
$ws.Range("B18").Value2 = 1
$ws.Range("A19").Value2 = "e"
$ws.Range("B19").Value2 = 163
$ws.Range("A20").Value2 = "s"
$ws.Range("B20").Value2 = 119
$ws.Range("A21").Value2 = "p"
$ws.Range("B21").Value2 = 33
$ws.Range("A22").Value2 = "q"
$ws.Range("B22").Value2 = 11
$ws.Range("A23").Value2 = "C"
$ws.Range("B23").Value2 = 1
$ws.Range("A24").Value2 = "F"
$ws.Range("B24").Value2 = 2
$ws.Range("A25").Value2 = "L"
$ws.Range("B25").Value2 = 2
$ws.Range("A26").Value2 = "m"
$ws.Range("B26").Value2 = 55
$ws.Range("A27").Value2 = "f"
$ws.Range("B27").Value2 = 10
$ws.Range("A28").Value2 = "t"
$ws.Range("B28").Value2 = 111
$ws.Range("A29").Value2 = "b"
$ws.Range("B29").Value2 = 18
$ws.Range("A30").Value2 = "S"
$ws.Range("B30").Value2 = 2
$ws.Range("A31").Value2 = "r"
$ws.Range("B31").Value2 = 76
$ws.Range("A32").Value2 = "a"
$ws.Range("B32").Value2 = 96
$ws.Range("A33").Value2 = "d"
$ws.Range("B33").Value2 = 51
$ws.Range("A34").Value2 = "M"
$ws.Range("B34").Value2 = 2
$ws.Range("A35").Value2 = "E"
$ws.Range("B35").Value2 = 1
$ws.Range("A36").Value2 = "N"
$ws.Range("B36").Value2 = 4
$ws.Range("A37").Value2 = "U"
$ws.Range("B37").Value2 = 1
$ws.Range("A38").Value2 = " "
$ws.Range("B38").Value2 = 245
$ws.Range("A39").Value2 = "u"
$ws.Range("B39").Value2 = 113
